$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.985.39"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.585.16"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.32"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.42"
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +2.69%  "

$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("E10").Value = "  +3.12%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.38"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").Value = "3.048.46"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "62.855.47"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("D17").Value = "2.584.98"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.31"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.11"
$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.70"
$ws.Range("E23").Value = "  -0.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.33"
$ws.Range("E24").Value = "  +3.30%  "

$ws.Range("D25").Value = "2.719.22"
$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.88"
$ws.Range("E29").Value = "  +8.99%  "

$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("E32").Value = "  +4.45%  "

$ws.Range("D33").Value = "0.0₃0823"
$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "466.99"
$ws.Range("E34").Value = "  +15.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.00"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("E36").Value = "  +4.81%  "

$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("E40").Value = "  +5.53%  "

$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "158.40"
$ws.Range("E43").Value = "  +5.01%  "

$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("E45").Value = "  +6.63%  "

$ws.Range("E46").Value = "  +2.82%  "

$ws.Range("E47").Value = "  +2.01%  "

$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.48"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +1.27%  "

Write-Host "Applied cryptos update"